# Update Name of Algo
# Applies updated RandomForest imputation results to columns C and D
# for the specific rows that changed between runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.169899999999997
$ws.Range("C3").Value = -12.3328
$ws.Range("C14").Value = -13.15419999999999
$ws.Range("C16").Value = -13.43530000000001
$ws.Range("D18").Value = -8.972699999999998
$ws.Range("C21").Value = -11.8979
$ws.Range("C23").Value = -12.6322
$ws.Range("D24").Value = -7.284899999999998
$ws.Range("C25").Value = -12.793
$ws.Range("D25").Value = -8.526600000000006
$ws.Range("C26").Value = -12.90070000000001
$ws.Range("D27").Value = -8.743100000000005
$ws.Range("C29").Value = -10.64830000000001
$ws.Range("D30").Value = -7.4824
$ws.Range("D31").Value = -8.483000000000006
$ws.Range("D39").Value = -8.183099999999998
$ws.Range("C40").Value = -13.4747
$ws.Range("D42").Value = -8.850099999999999
$ws.Range("D48").Value = -7.389799999999999
$ws.Range("D51").Value = -7.779999999999998
$ws.Range("D52").Value = -7.805099999999999
$ws.Range("C53").Value = -10.30450000000001
$ws.Range("D55").Value = -8.835900000000001
$ws.Range("D56").Value = -7.8531
$ws.Range("C57").Value = -14.15199999999999
$ws.Range("D57").Value = -8.209199999999996
$ws.Range("C59").Value = -12.78149999999999
$ws.Range("D60").Value = -8.101400000000002
$ws.Range("C65").Value = -12.10360000000001
$ws.Range("C69").Value = -10.7239
$ws.Range("D73").Value = -7.558800000000002
$ws.Range("D74").Value = -8.539600000000005
$ws.Range("C79").Value = -10.72330000000002
$ws.Range("C83").Value = -14.1185
$ws.Range("D89").Value = -6.033999999999999
$ws.Range("D90").Value = -8.081200000000001
$ws.Range("C91").Value = -10.173
$ws.Range("D92").Value = -6.425300000000001
$ws.Range("C93").Value = -11.61260000000001
$ws.Range("C100").Value = -13.56169999999998
